$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data for REMIND-MAgPIE 2.1-4.2
# (force text format on A6 first so the dd.mm.yyyy-looking date string is
# stored as literal text, not auto-converted to a date serial number)
$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "10.06.2020"
$ws.Range("A6").NumberFormat = "General"
$ws.Range("B6").Value = "REMIND-MAgPIE 2.1-4.2"
$ws.Range("C6").Value = "IPCC_AR6_model_registration_REMIND-MAgPIE_2.1-4.2.xlsx"
$ws.Range("D6").Value = "Bjoern Soergel, Alois Dirnaichner, Isabelle Weindl"
$ws.Range("E6").Value = "Coupled system with latest REMIND and MAgPIE versions (e.g. for SDP)"

# Match row 5's formatting: column E wraps text, like the rest of the table
$ws.Range("E6").WrapText = $true

# Adjust column widths to match new content
# (target widths from the source file are 56.45 / 45.46 / 73.27 chars; the
# COM width model here snaps to the nearest 1/6-character pixel boundary,
# so these inputs are chosen to land on the closest achievable width)
$ws.Columns.Item(3).ColumnWidth = 55.66666666666667
$ws.Columns.Item(4).ColumnWidth = 44.66666666666667
$ws.Columns.Item(5).ColumnWidth = 72.5

# Set the active cell selection to C6
$ws.Range("C6").Select()
